$wb = $excel.ActiveWorkbook

# --- Summary ---
$ws = $wb.Worksheets.Item("Summary")
$ws.Range("B4").Value = "inf"
$ws.Range("B6").Value = 72148.09464184758
$ws.Range("B7").Value = 10208330.5779121
$ws.Range("B8").Value = 24483963.20213782
$ws.Range("B10").Value = 3060021.080760467

# --- Fed-in Capacity ---
$ws = $wb.Worksheets.Item("Fed-in Capacity")
$ws.Range("P18").Value = 65.92768427608706
$ws.Range("P20").Value = 135.4597561231036
$ws.Range("I21").Value = 0
$ws.Range("J21").Value = 0
$ws.Range("L21").Value = 0
$ws.Range("Q21").Value = 0
$ws.Range("L22").Value = 90.4687457914608
$ws.Range("M22").Value = 0
$ws.Range("N22").Value = 0
$ws.Range("P22").Value = 101.5955875616828
$ws.Range("K24").Value = 80.29914934735042
$ws.Range("L24").Value = 61.18167021676314
$ws.Range("M24").Value = 51.84373129681028
$ws.Range("O24").Value = 57.81213424001893
$ws.Range("L25").Value = 0
$ws.Range("L26").Value = 0
$ws.Range("M26").Value = 113.4004983079896
$ws.Range("Q26").Value = 150.3839754851235
$ws.Range("K27").Value = 80.29914934735042
$ws.Range("M27").Value = 51.84373129681028
$ws.Range("N27").Value = 38.66169381481656
$ws.Range("O27").Value = 57.81213424001893
$ws.Range("P27").Value = 0
$ws.Range("M28").Value = 92.09541281912071
$ws.Range("O28").Value = 96.22962838366004
$ws.Range("P28").Value = 101.5955875616828
$ws.Range("M30").Value = 0
$ws.Range("J32").Value = 124.5190384721106
$ws.Range("Q32").Value = 150.3839754851235
$ws.Range("J33").Value = 93.17061249236157
$ws.Range("L33").Value = 61.18167021676314
$ws.Range("M33").Value = 51.84373129681028
$ws.Range("N33").Value = 38.66169381481656
$ws.Range("P33").Value = 65.92768427608706
$ws.Range("Q33").Value = 94.49434172313325
$ws.Range("L36").Value = 0
$ws.Range("P36").Value = 0
$ws.Range("N38").Value = 0
$ws.Range("P41").Value = 135.4597561231036
$ws.Range("J42").Value = 0
$ws.Range("Q42").Value = 94.49434172313325
$ws.Range("R44").Value = 65.71641987298243
$ws.Range("R45").Value = 0
$ws.Range("K46").Value = 0

# --- Unmet Demand ---
$ws = $wb.Worksheets.Item("Unmet Demand")
$ws.Range("P18").Value = 0
$ws.Range("P20").Value = 0
$ws.Range("I21").Value = 87.25340171355576
$ws.Range("J21").Value = 93.17061249236157
$ws.Range("L21").Value = 61.18167021676314
$ws.Range("Q21").Value = 94.49434172313325
$ws.Range("L22").Value = 0
$ws.Range("M22").Value = 92.09541281912071
$ws.Range("N22").Value = 81.96869489115805
$ws.Range("P22").Value = 0
$ws.Range("K24").Value = 0
$ws.Range("L24").Value = 0
$ws.Range("M24").Value = 0
$ws.Range("O24").Value = 0
$ws.Range("L25").Value = 90.4687457914608
$ws.Range("L26").Value = 130.6648563030561
$ws.Range("M26").Value = 0
$ws.Range("Q26").Value = 0
$ws.Range("K27").Value = 0
$ws.Range("M27").Value = 0
$ws.Range("N27").Value = 0
$ws.Range("O27").Value = 0
$ws.Range("P27").Value = 65.92768427608706
$ws.Range("M28").Value = 0
$ws.Range("O28").Value = 0
$ws.Range("P28").Value = 0
$ws.Range("M30").Value = 51.84373129681028
$ws.Range("J32").Value = 0
$ws.Range("Q32").Value = 0
$ws.Range("J33").Value = 0
$ws.Range("L33").Value = 0
$ws.Range("M33").Value = 0
$ws.Range("N33").Value = 0
$ws.Range("P33").Value = 0
$ws.Range("Q33").Value = 0
$ws.Range("L36").Value = 61.18167021676314
$ws.Range("P36").Value = 65.92768427608706
$ws.Range("N38").Value = 110.5750244233121
$ws.Range("P41").Value = 0
$ws.Range("J42").Value = 93.17061249236157
$ws.Range("Q42").Value = 0
$ws.Range("R44").Value = 108.0327934026353
$ws.Range("R45").Value = 123.5547069419379
$ws.Range("K46").Value = 94.30397654773019

# --- Household Surplus ---
$ws = $wb.Worksheets.Item("Household Surplus")
$ws.Range("B7").Value = 192657.4251299174
$ws.Range("B8").Value = 171531.8174962914
$ws.Range("B9").Value = 198637.9880329445
$ws.Range("B10").Value = 309745.1289038616
$ws.Range("B11").Value = 279503.3814941623
$ws.Range("B12").Value = 278856.3690631902
$ws.Range("B13").Value = 180607.232914963
$ws.Range("B14").Value = 221044.0398916555
$ws.Range("B15").Value = 204796.4608921716
$ws.Range("B16").Value = 154301.354011648

# --- Costs and Revenues ---
$ws = $wb.Worksheets.Item("Costs and Revenues")
$ws.Range("G2").Value = 65343.16542154989
$ws.Range("H2").Value = 59609.07192099428
$ws.Range("I2").Value = 66966.46106665728
$ws.Range("J2").Value = 97124.1135887634
$ws.Range("K2").Value = 88915.63929184499
$ws.Range("L2").Value = 88740.02163200968
$ws.Range("M2").Value = 62072.39896320511
$ws.Range("N2").Value = 73048.10371402171
$ws.Range("O2").Value = 68638.04655701894
$ws.Range("P2").Value = 54932.23183230536
$ws.Range("E3").Value = 133100.0000000001
$ws.Range("G4").Value = 17262.61268259926
$ws.Range("H4").Value = 11528.51918204364
$ws.Range("I4").Value = 18885.90832770665
$ws.Range("J4").Value = 49043.56084981273
$ws.Range("K4").Value = 40835.08655289431
$ws.Range("L4").Value = 40659.46889305903
$ws.Range("M4").Value = 13991.84622425448
$ws.Range("N4").Value = 24967.55097507107
$ws.Range("O4").Value = 20557.49381806829
$ws.Range("P4").Value = 6851.679093354725
$ws.Range("E6").Value = -88631.0593660577
$ws.Range("J6").Value = 44468.94063394235
$ws.Range("K6").Value = 44468.94063394235
$ws.Range("L6").Value = 44468.94063394232
$ws.Range("N6").Value = 44468.94063394232
$ws.Range("O6").Value = 44468.94063394232
